$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.734.39"
$ws.Range("E2").Value = "  +3.54%  "

$ws.Range("D3").Value = "2.253.34"
$ws.Range("E3").Value = "  +2.90%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'253.81"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("E6").Value = "  +1.38%  "

$ws.Range("D7").Value = "'71.89"
$ws.Range("E7").Value = "  +4.98%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.648"
$ws.Range("E9").Value = "  +12.75%  "

$ws.Range("D10").Value = "'41.12"
$ws.Range("E10").Value = "  +9.23%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0967"
$ws.Range("E11").Value = "  +3.12%  "

$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").Value = "'59.26"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").Value = "'7.42"
$ws.Range("E13").Value = "  +4.29%  "

$ws.Range("D14").Value = "'0.104"

$ws.Range("D15").Value = "2.592.54"
$ws.Range("E15").Value = "  +3.03%  "

$ws.Range("D16").Value = "'0.888"
$ws.Range("E16").Value = "  +1.31%  "

$ws.Range("D17").Value = "'14.79"
$ws.Range("E17").Value = "  +2.03%  "

$ws.Range("D18").Value = "2.250.53"
$ws.Range("E18").Value = "  +3.53%  "

$ws.Range("D19").Value = "42.724.33"
$ws.Range("E19").Value = "  +3.42%  "

$ws.Range("E20").Value = "  +1.82%  "

$ws.Range("E21").Value = "  +0.99%  "

$ws.Range("D22").Value = "'73.24"
$ws.Range("E22").Value = "  +1.61%  "

$ws.Range("D23").Value = "'235.03"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("D24").Value = "'2.10"
$ws.Range("E24").Value = "  +3.95%  "

$ws.Range("D25").Value = "'3.96"
$ws.Range("E25").Value = "  +0.63%  "

$ws.Range("D26").Value = "'11.69"
$ws.Range("E26").Value = "  -1.11%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("E28").Value = "  -2.73%  "

$ws.Range("D29").Value = "'3.69"
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +2.22%  "

$ws.Range("D31").Value = "'167.82"
$ws.Range("E31").Value = "  -0.76%  "

$ws.Range("D32").Value = "'21.03"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("E33").Value = "  +9.63%  "

$ws.Range("D34").Value = "'6.18"
$ws.Range("E34").Value = "  +12.40%  "

$ws.Range("D35").Value = "'0.0786"
$ws.Range("E35").Value = "  +3.80%  "

$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").Value = "'28.89"
$ws.Range("E37").Value = "  +8.20%  "

$ws.Range("D38").Value = "'4.71"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("E39").Value = "  -1.15%  "

$ws.Range("E40").Value = "  +7.77%  "

$ws.Range("E41").Value = "  +4.08%  "

$ws.Range("D42").Value = "'5.85"
$ws.Range("E42").Value = "  +2.94%  "

$ws.Range("D43").Value = "'12.48"
$ws.Range("E43").Value = "  +0.82%  "

$ws.Range("D44").Value = "'64.21"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.202"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.03"
$ws.Range("E47").Value = "  +4.18%  "

$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("E51").Value = "  +1.20%  "
